# Apply the "4 mdelo melhores rstds" edit:
# - Column A (model name) values for data rows 2..26 are re-ordered (the models
#   were re-ranked by rstd, so the labels attached to each row change).
# - Columns B..I for every data row (2..26) are overwritten with the metrics
#   of the 4 best models (duplicated across all rows, matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-name ordering for rows 2..26
$names = @{
    2  = "model_2_2_0"
    3  = "model_2_2_22"
    4  = "model_2_2_21"
    5  = "model_2_2_20"
    6  = "model_2_2_19"
    7  = "model_2_2_18"
    8  = "model_2_2_17"
    9  = "model_2_2_16"
    10 = "model_2_2_15"
    11 = "model_2_2_14"
    12 = "model_2_2_13"
    13 = "model_2_2_23"
    14 = "model_2_2_12"
    15 = "model_2_2_10"
    16 = "model_2_2_9"
    17 = "model_2_2_8"
    18 = "model_2_2_7"
    19 = "model_2_2_6"
    20 = "model_2_2_5"
    21 = "model_2_2_4"
    22 = "model_2_2_3"
    23 = "model_2_2_2"
    24 = "model_2_2_1"
    25 = "model_2_2_11"
    26 = "model_2_2_24"
}

# New metric values shared by every data row (2..26)
$b = -0.001236788884735551
$c = -0.9459796395945301
$d = -0.1472438604251296
$e = -0.09099460614211563
$f = 1.108074069023132
$g = 0.7651928067207336
$h = 2.127733945846558
$i = 1.406388759613037

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("A$row").Value = $names[$row]
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
    $ws.Range("I$row").Value = $i
}
